$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Text fix-ups: drop the trailing colon from the inbuilt "Example" labels
#    (the "Exercise" labels keep their colon, they are untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Example 1.1:", $false, $false, $false, $false, $false, $true, 1, $false, "Example 1.1", 2) | Out-Null
$d.Content.Find.Execute("Example 2.1:", $false, $false, $false, $false, $false, $true, 1, $false, "Example 2.1", 2) | Out-Null
$d.Content.Find.Execute("Example 2.2:", $false, $false, $false, $false, $false, $true, 1, $false, "Example 2.2", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Re-introduce the bookmarks that mark each user-defined "Example"/
#    "Exercise" environment individually (these were previously lost, per
#    the commit message, so links to them were broken).
#
#    Paragraph numbering is stable across these calls because adding a
#    bookmark never inserts or removes a paragraph.
# ---------------------------------------------------------------------------

# -- Section 1 --------------------------------------------------------------
# Example 1.1 (first, builtin-numbered example)
$p = $d.Paragraphs.Item(8)
$rng = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("exm:inbuilt", $rng) | Out-Null

# Exercise 1.1 (author supplied id) spans the name paragraph + the body paragraph
$pStart = $d.Paragraphs.Item(11)
$pEnd = $d.Paragraphs.Item(12)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$d.Bookmarks.Add("Exe:author", $rng) | Out-Null

# Example 1.1 (second, builtin-numbered example, repeated on purpose)
$p = $d.Paragraphs.Item(14)
$rng = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("exm:inbuilt", $rng) | Out-Null

# -- Section 2 ----------------------------------------------------------------
# Example 2.1
$p = $d.Paragraphs.Item(18)
$rng = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("exm:unnamed-chunk-1", $rng) | Out-Null

# Example 1.1 (repeated, builtin-numbered example)
$p = $d.Paragraphs.Item(21)
$rng = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("exm:inbuilt", $rng) | Out-Null

# Exercise 1.1 (author supplied id)
$pStart = $d.Paragraphs.Item(24)
$pEnd = $d.Paragraphs.Item(25)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$d.Bookmarks.Add("Exe:author", $rng) | Out-Null

# Example 2.2
$p = $d.Paragraphs.Item(27)
$rng = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("exm:unnamed-chunk-2", $rng) | Out-Null

# Exercise 2.1 (author supplied id, second author label)
$pStart = $d.Paragraphs.Item(30)
$pEnd = $d.Paragraphs.Item(31)
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$d.Bookmarks.Add("Exe:author2", $rng) | Out-Null
